# words_dict.xlsx - "update framework: new functions"
#
# 1. New translation row for the template auto-update feature
#    (TEMPLATE_AUTO_UPDATE / 自动更新 / Auto Update), inserted at row 163
#    (just above SAMPLING), pushing all following rows down by one.
# 2. The EXCLUSIVE_WARNING message text is updated from "server mode" to
#    "distributed mode" (Chinese + English) - it now sits at row 402 because
#    of the insertion above.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- 1. insert the new TEMPLATE_AUTO_UPDATE row at 163 -----------------
$ws.Rows("163:163").Insert()
$ws.Range("A163").Value = "TEMPLATE_AUTO_UPDATE"
$ws.Range("B163").Value = "自动更新"
$ws.Range("C163").Value = "Auto Update"

# --- 2. EXCLUSIVE_WARNING text: server mode -> distributed mode --------
# (row 401 shifted to row 402 by the insertion above)
$ws.Range("B402").Value = "您当前运行于分布模式，独占模式可能造成MATLAB崩溃，继续？"
$ws.Range("C402").Value = "You are currently running in distributed mode. Exclusive mode may cause MATLAB to crash. Continue?"

# --- 3. cosmetic: restore the view/selection state recorded in the file ---
$excel.ActiveWindow.ScrollRow = 259
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A279").Select()
